# Rename the three worksheets (spaces removed from the Spanish names,
# translated to English short titles). Excel automatically rewrites any
# formulas / defined names that reference the old sheet names.
$wb = $excel.ActiveWorkbook

$wsConnectivity = $wb.Worksheets.Item("CONEXIONES CON")
$wsCoordinates  = $wb.Worksheets.Item("COORDENADAS COOR")
$wsFreeNodes    = $wb.Worksheets.Item("NODOS LIBRES NL")

$wsCoordinates.Name  = "COORDINATES"
$wsConnectivity.Name = "CONNECTIVITY"
$wsFreeNodes.Name    = "FREE NODES"

# Move the active tab from "FREE NODES" (previously NODOS LIBRES NL) to
# "CONNECTIVITY" (previously CONEXIONES CON).
$wsConnectivity.Activate()
